$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B2").Value = 0.9351974020079297
$ws.Range("C2").Value = 0.06415970883649891
$ws.Range("E2").Value = 0.04764472388581886
$ws.Range("F2").Value = 0.4443680307746263
$ws.Range("G2").Value = 0.002603767831384336
$ws.Range("I2").Value = 2.920195067636683
$ws.Range("K2").Value = 0.8281819962446946
$ws.Range("L2").Value = 0.2641690103052809
$ws.Range("B3").Value = 0.9152238243558486
$ws.Range("C3").Value = 0.05721839044161925
$ws.Range("E3").Value = 0.04740006353959458
$ws.Range("F3").Value = 0.387822817061874
$ws.Range("G3").Value = 0.002608239461068854
$ws.Range("I3").Value = 2.845266988450732
$ws.Range("K3").Value = 0.799705497978465
$ws.Range("L3").Value = 0.2568509453645191
$ws.Range("B4").Value = 0.9037011543522624
$ws.Range("C4").Value = 0.05297132539777749
$ws.Range("E4").Value = 0.04727484433793983
$ws.Range("F4").Value = 0.3531389305168915
$ws.Range("G4").Value = 0.002611127471306838
$ws.Range("I4").Value = 2.79955725436497
$ws.Range("K4").Value = 0.7828389731793948
$ws.Range("L4").Value = 0.2525065330281109
$ws.Range("B5").Value = 0.8991918857648784
$ws.Range("C5").Value = 0.05124423535669109
$ws.Range("E5").Value = 0.04723011184620951
$ws.Range("F5").Value = 0.3390132514313251
$ws.Range("G5").Value = 0.002612340289612146
$ws.Range("I5").Value = 2.781003120320293
$ws.Range("K5").Value = 0.7761208958891928
$ws.Range("L5").Value = 0.2507735291021049
$ws.Range("B6").Value = 0.8984543763788793
$ws.Range("C6").Value = 0.05095766839805549
$ws.Range("E6").Value = 0.04722306451628278
$ws.Range("F6").Value = 0.336668177824194
$ws.Range("G6").Value = 0.002612543850818132
$ws.Range("I6").Value = 2.777926580494992
$ws.Range("K6").Value = 0.7750147283492481
$ws.Range("L6").Value = 0.2504880208929592
$ws.Range("B7").Value = 0.9036395864869746
$ws.Range("C7").Value = 0.05294801875785993
$ws.Range("E7").Value = 0.04727421555934441
$ws.Range("F7").Value = 0.3529483938344953
$ws.Range("G7").Value = 0.002611143682122653
$ws.Range("I7").Value = 2.799306733546416
$ws.Range("K7").Value = 0.7827477429244709
$ws.Range("L7").Value = 0.2524830098663955
$ws.Range("B8").Value = 0.9281566267720507
$ws.Range("C8").Value = 0.06176316165584694
$ws.Range("E8").Value = 0.047555180150475
$ws.Range("F8").Value = 0.4248636149813336
$ws.Range("G8").Value = 0.00260528016433148
$ws.Range("I8").Value = 2.894297373770925
$ws.Range("K8").Value = 0.8182348214706678
$ws.Range("L8").Value = 0.2616147928085866
$ws.Range("B9").Value = 0.9821232804170563
$ws.Range("C9").Value = 0.07917489860943761
$ws.Range("E9").Value = 0.04830428710752699
$ws.Range("F9").Value = 0.5661985755041457
$ws.Range("G9").Value = 0.002594906293427002
$ws.Range("I9").Value = 3.083005725886309
$ws.Range("K9").Value = 0.8927474928094057
$ws.Range("L9").Value = 0.2807079925585612
$ws.Range("B10").Value = 1.025380489954131
$ws.Range("C10").Value = 0.09205429731829895
$ws.Range("E10").Value = 0.04897525586233087
$ws.Range("F10").Value = 0.6702781546542269
$ws.Range("G10").Value = 0.002587962391481025
$ws.Range("I10").Value = 3.223255973192295
$ws.Range("K10").Value = 0.9505254722756717
$ws.Range("L10").Value = 0.2954665166327572
$ws.Range("B11").Value = 1.045847248287259
$ws.Range("C11").Value = 0.09793473507701833
$ws.Range("E11").Value = 0.04930665799527034
$ws.Range("F11").Value = 0.7176906081379002
$ws.Range("G11").Value = 0.002584948954682992
$ws.Range("I11").Value = 3.28743552386598
$ws.Range("K11").Value = 0.9774765621754682
$ws.Range("L11").Value = 0.3023410705019387
$ws.Range("B12").Value = 1.053711165303099
$ws.Range("C12").Value = 0.1001647673293746
$ws.Range("E12").Value = 0.04943591077497089
$ws.Range("F12").Value = 0.7356546913071611
$ws.Range("G12").Value = 0.002583828624675429
$ws.Range("I12").Value = 3.31179506207809
$ws.Range("K12").Value = 0.9877787264607321
$ws.Range("L12").Value = 0.3049675257259139
$ws.Range("B13").Value = 1.052012474900437
$ws.Range("C13").Value = 0.09968434370225054
$ws.Range("E13").Value = 0.04940790683127005
$ws.Range("F13").Value = 0.7317853510981394
$ws.Range("G13").Value = 0.002584068984893876
$ws.Range("I13").Value = 3.306546277005509
$ws.Range("K13").Value = 0.9855556788450599
$ws.Range("L13").Value = 0.3044008379264938
$ws.Range("B14").Value = 1.046491939764991
$ws.Range("C14").Value = 0.09811813546127723
$ws.Range("E14").Value = 0.04931721641692022
$ws.Range("F14").Value = 0.7191683204515869
$ws.Range("G14").Value = 0.002584856368452551
$ws.Range("I14").Value = 3.289438463336978
$ws.Range("K14").Value = 0.9783221938274664
$ws.Range("L14").Value = 0.3025566848834131
$ws.Range("B15").Value = 1.043125251797505
$ws.Range("C15").Value = 0.09715921360674429
$ws.Range("E15").Value = 0.04926215514758781
$ws.Range("F15").Value = 0.7114413442032514
$ws.Range("G15").Value = 0.002585341368300531
$ws.Range("I15").Value = 3.278966792473284
$ws.Range("K15").Value = 0.9739040370583325
$ws.Range("L15").Value = 0.3014301131729979
$ws.Range("B16").Value = 1.024058826468064
$ws.Range("C16").Value = 0.09167044227430665
$ws.Range("E16").Value = 0.04895412414775535
$ws.Range("F16").Value = 0.6671810134426437
$ws.Range("G16").Value = 0.00258816224225104
$ws.Range("I16").Value = 3.219069454924067
$ws.Range("K16").Value = 0.9487776204133809
$ws.Range("L16").Value = 0.2950204917620027
$ws.Range("B17").Value = 1.012564346058923
$ws.Range("C17").Value = 0.08830887461240877
$ws.Range("E17").Value = 0.04877185685170282
$ws.Range("F17").Value = 0.6400460337215605
$ws.Range("G17").Value = 0.002589929911825628
$ws.Range("I17").Value = 3.182422614231882
$ws.Range("K17").Value = 0.9335346257039134
$ws.Range("L17").Value = 0.2911296474521379
$ws.Range("B18").Value = 1.006027263257778
$ws.Range("C18").Value = 0.08637740821345119
$ws.Range("E18").Value = 0.04866948558730755
$ws.Range("F18").Value = 0.6244449056556647
$ws.Range("G18").Value = 0.002590960319401547
$ws.Range("I18").Value = 3.161379897189846
$ws.Range("K18").Value = 0.9248300726944194
$ws.Range("L18").Value = 0.2889068693644532
$ws.Range("B19").Value = 1.003826664995245
$ws.Range("C19").Value = 0.08572378925455837
$ws.Range("E19").Value = 0.04863524782725115
$ws.Range("F19").Value = 0.6191636801734006
$ws.Range("G19").Value = 0.002591311552397188
$ws.Range("I19").Value = 3.154261255277007
$ws.Range("K19").Value = 0.9218936379227216
$ws.Range("L19").Value = 0.2881568714191047
$ws.Range("B20").Value = 1.013780268591262
$ws.Range("C20").Value = 0.08866650930835362
$ws.Range("E20").Value = 0.04879100453911001
$ws.Range("F20").Value = 0.642933953830422
$ws.Range("G20").Value = 0.002589740324287675
$ws.Range("I20").Value = 3.186320039140412
$ws.Range("K20").Value = 0.9351507648188431
$ws.Range("L20").Value = 0.2915422676960731
$ws.Range("B21").Value = 1.048110370297763
$ws.Range("C21").Value = 0.09857807977203947
$ws.Range("E21").Value = 0.04934375244559419
$ws.Range("F21").Value = 0.7228739723492197
$ws.Range("G21").Value = 0.002584624531345433
$ws.Range("I21").Value = 3.294461906385521
$ws.Range("K21").Value = 0.9804442251976866
$ws.Range("L21").Value = 0.3030977267167145
$ws.Range("B22").Value = 1.071209309569724
$ws.Range("C22").Value = 0.1050748047310606
$ws.Range("E22").Value = 0.04972690697739068
$ws.Range("F22").Value = 0.7751780083420101
$ws.Range("G22").Value = 0.002581402209721376
$ws.Range("I22").Value = 3.365466826753789
$ws.Range("K22").Value = 1.010608028408427
$ws.Range("L22").Value = 0.3107852186347344
$ws.Range("B23").Value = 1.058820318233842
$ws.Range("C23").Value = 0.1016056001986101
$ws.Range("E23").Value = 0.04952040807810221
$ws.Range("F23").Value = 0.7472568307915566
$ws.Range("G23").Value = 0.002583110975591192
$ws.Range("I23").Value = 3.327539615886025
$ws.Range("K23").Value = 0.9944574980097514
$ws.Range("L23").Value = 0.3066698482741401
$ws.Range("B24").Value = 1.013230328100889
$ws.Range("C24").Value = 0.08850481911127872
$ws.Range("E24").Value = 0.04878234034201867
$ws.Range("F24").Value = 0.6416283278902313
$ws.Range("G24").Value = 0.002589825992906662
$ws.Range("I24").Value = 3.184557932345399
$ws.Range("K24").Value = 0.9344199250514862
$ws.Range("L24").Value = 0.2913556781101363
$ws.Range("B25").Value = 0.9668916270301509
$ws.Range("C25").Value = 0.07445005740525801
$ws.Range("E25").Value = 0.0480804424664143
$ws.Range("F25").Value = 0.5279251897347308
$ws.Range("G25").Value = 0.00259759311918759
$ws.Range("I25").Value = 3.031681977967295
$ws.Range("K25").Value = 0.8720595218709946
$ws.Range("L25").Value = 0.2754150244351479
